$wb = $excel.ActiveWorkbook

# Add the new "token" worksheet by duplicating the existing "Template" sheet
# (same layout/columns/styles), then renaming the copy.
$template = $wb.Worksheets.Item("Template")
$template.Copy($null, $template)
$tokenSheet = $wb.Worksheets.Item("Template (2)")
$tokenSheet.Name = "token"

# Title cell
$tokenSheet.Range("A1").Value = "token"

# Fill in the "Type" column first (all rows use the same value)
$tokenSheet.Range("E5").Value = "varchar()"
$tokenSheet.Range("E6").Value = "varchar()"
$tokenSheet.Range("E7").Value = "varchar()"

# Logic name for the token column
$tokenSheet.Range("D5").Value = "Token"

# Column names
$tokenSheet.Range("C5").Value = "token"
$tokenSheet.Range("C6").Value = "id"
$tokenSheet.Range("C7").Value = "role"

# Remaining logic names
$tokenSheet.Range("D6").Value = "连接其他表"
$tokenSheet.Range("D7").Value = "连接表名"

# Numeric / flag columns
$tokenSheet.Range("B5").Value = 0
$tokenSheet.Range("F5").Value = 1
$tokenSheet.Range("I5").Value = 1

$tokenSheet.Range("B6").Value = 1
$tokenSheet.Range("I6").Value = 1

$tokenSheet.Range("B7").Value = 2
$tokenSheet.Range("I7").Value = 1

# Update view/selection state: Template is no longer the active tab and has
# its whole sheet selected; "token" becomes the active tab with M8 selected.
$template.Activate()
$template.Cells.Select()

$tokenSheet.Activate()
$tokenSheet.Range("M8").Select()
